$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidateFormulas")

# Insert a new row at 31, shifting existing rows 31-41 down to 32-42
$ws.Range("A31").EntireRow.Insert()

# Populate the new row 31
$ws.Range("A31").Value = "If"
$ws.Range("A31").Font.Bold = $true

$ws.Range("B31").Formula = "=IF(B2>3,B3,B5)"
$ws.Range("C31").Formula = "=IF((B2*B3)*C1<0,(B2*B3)*C1,ABS((B2*B3)*C1))"
$ws.Range("D31").Formula = "=IF((B2*B3)*C1<0,ABS((B2*B3)*C1),(B2*B3)*C1)"

# Update selection to match target (D31 active cell)
$ws.Range("D31").Select()
